$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "PERCENT COMPLETE" helper column (E) that flags each row as
# complete (1) or not (0) based on the boolean COMPLETE column (D).
$ws.Range("E2").Formula = "=IF(D2,1,0)"
$ws.Range("E3:E55").Formula = "=IF(D3,1,0)"

# Add summary row below the data: a right-aligned label and the
# percent-complete calculation averaged across all the rows.
$ws.Range("C57").Value = "PERCENT COMPLETE"
$ws.Range("C57").HorizontalAlignment = -4152
$ws.Range("E57").Formula = "=100*AVERAGE(E2:E55)"

# Match the workbook's last active selection.
$ws.Range("E53").Select()
